$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 9753
$ws.Range("E2").Value = 970
$ws.Range("F2").Value = 970
$ws.Range("G2").Value = 1167
$ws.Range("H2").Value = 868
$ws.Range("I2").Value = 839
$ws.Range("J2").Value = 28
$ws.Range("K2").Value = 13304
$ws.Range("L2").Value = 4097
$ws.Range("M2").Value = 9206
$ws.Range("N2").Value = 8927
$ws.Range("O2").Value = 280
$ws.Range("P2").Value = 584
$ws.Range("Q2").Value = 568
$ws.Range("R2").Value = -1634
$ws.Range("S2").Value = 1094
$ws.Range("T2").Value = 585
$ws.Range("U2").Value = -17
$ws.Range("V2").Value = 1605
$ws.Range("W2").Value = 9.94
$ws.Range("X2").Value = 8.9
$ws.Range("Y2").Value = 9.789999999999999
$ws.Range("Z2").Value = 7.13
$ws.Range("AA2").Value = 44.51
$ws.Range("AB2").Value = 1420.34
$ws.Range("AC2").Value = 7183
$ws.Range("AD2").Value = 19.07
$ws.Range("AE2").Value = 77199
$ws.Range("AF2").Value = 1.77
$ws.Range("AG2").Value = 1250
$ws.Range("AH2").Value = 0.91
$ws.Range("AI2").Value = 17.22
$ws.Range("AJ2").Value = 11686538

$ws.Range("D3").Value = 10478
$ws.Range("E3").Value = 917
$ws.Range("F3").Value = 917
$ws.Range("G3").Value = 1289
$ws.Range("H3").Value = 957
$ws.Range("I3").Value = 950
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 14151
$ws.Range("L3").Value = 4174
$ws.Range("M3").Value = 9977
$ws.Range("N3").Value = 9568
$ws.Range("O3").Value = 409
$ws.Range("P3").Value = 584
$ws.Range("Q3").Value = 452
$ws.Range("R3").Value = 78
$ws.Range("S3").Value = -276
$ws.Range("T3").Value = 1212
$ws.Range("U3").Value = -759
$ws.Range("V3").Value = 1327
$ws.Range("W3").Value = 8.75
$ws.Range("X3").Value = 9.130000000000001
$ws.Range("Y3").Value = 10.27
$ws.Range("Z3").Value = 6.97
$ws.Range("AA3").Value = 41.84
$ws.Range("AB3").Value = 1551.82
$ws.Range("AC3").Value = 8126
$ws.Range("AD3").Value = 22.52
$ws.Range("AE3").Value = 82746
$ws.Range("AF3").Value = 2.21
$ws.Range("AG3").Value = 1750
$ws.Range("AH3").Value = 0.96
$ws.Range("AI3").Value = 21.31
$ws.Range("AJ3").Value = 11686538

$ws.Range("D4").Value = 11979
$ws.Range("E4").Value = 785
$ws.Range("F4").Value = 785
$ws.Range("G4").Value = 795
$ws.Range("H4").Value = 652
$ws.Range("I4").Value = 630
$ws.Range("J4").Value = 22
$ws.Range("K4").Value = 15506
$ws.Range("L4").Value = 4916
$ws.Range("M4").Value = 10590
$ws.Range("N4").Value = 9931
$ws.Range("O4").Value = 659
$ws.Range("P4").Value = 584
$ws.Range("Q4").Value = -17
$ws.Range("R4").Value = -1111
$ws.Range("S4").Value = 1106
$ws.Range("T4").Value = 807
$ws.Range("U4").Value = -824
$ws.Range("V4").Value = 2488
$ws.Range("W4").Value = 6.55
$ws.Range("X4").Value = 5.44
$ws.Range("Y4").Value = 6.46
$ws.Range("Z4").Value = 4.39
$ws.Range("AA4").Value = 46.42
$ws.Range("AB4").Value = 1645.95
$ws.Range("AC4").Value = 5388
$ws.Range("AD4").Value = 29.14
$ws.Range("AE4").Value = 87011
$ws.Range("AF4").Value = 1.8
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 0.8
$ws.Range("AI4").Value = 22.66
$ws.Range("AJ4").Value = 11686538

$ws.Range("D5").Value = 12879
$ws.Range("E5").Value = 903
$ws.Range("F5").Value = 903
$ws.Range("G5").Value = 722
$ws.Range("H5").Value = 567
$ws.Range("I5").Value = 532
$ws.Range("J5").Value = 35
$ws.Range("K5").Value = 17036
$ws.Range("L5").Value = 5949
$ws.Range("M5").Value = 11087
$ws.Range("N5").Value = 10402
$ws.Range("O5").Value = 685
$ws.Range("P5").Value = 584
$ws.Range("Q5").Value = 579
$ws.Range("R5").Value = 148
$ws.Range("S5").Value = 482
$ws.Range("T5").Value = 295
$ws.Range("U5").Value = 284
$ws.Range("V5").Value = 3130
$ws.Range("W5").Value = 7.01
$ws.Range("X5").Value = 4.4
$ws.Range("Y5").Value = 5.24
$ws.Range("Z5").Value = 3.48
$ws.Range("AA5").Value = 53.65
$ws.Range("AB5").Value = 1714.5
$ws.Range("AC5").Value = 4556
$ws.Range("AD5").Value = 49.6
$ws.Range("AE5").Value = 91141
$ws.Range("AF5").Value = 2.48
$ws.Range("AG5").Value = 1250
$ws.Range("AH5").Value = 0.55
$ws.Range("AI5").Value = 26.79
$ws.Range("AJ5").Value = 11686538

$ws.Range("D6").Value = 13349
$ws.Range("E6").Value = 502
$ws.Range("F6").Value = 502
$ws.Range("G6").Value = 502
$ws.Range("H6").Value = 342
$ws.Range("I6").Value = 344
$ws.Range("K6").Value = 17234
$ws.Range("L6").Value = 5972
$ws.Range("M6").Value = 11262
$ws.Range("N6").Value = 10596
$ws.Range("P6").Value = 584
$ws.Range("Q6").Value = 144
$ws.Range("R6").Value = -861
$ws.Range("S6").Value = -45
$ws.Range("T6").Value = 414
$ws.Range("U6").Value = -270
$ws.Range("V6").Value = 3288
$ws.Range("W6").Value = 3.76
$ws.Range("X6").Value = 2.57
$ws.Range("Y6").Value = 3.28
$ws.Range("Z6").Value = 2
$ws.Range("AA6").Value = 53.03
$ws.Range("AB6").Value = 1760.95
$ws.Range("AC6").Value = 2946
$ws.Range("AD6").Value = 46.16
$ws.Range("AE6").Value = 92841
$ws.Range("AF6").Value = 1.46
$ws.Range("AI6").Value = 33.15
$ws.Range("AJ6").Value = 11686538

$ws.Range("D7").Value = 13636
$ws.Range("E7").Value = 572
$ws.Range("G7").Value = 258
$ws.Range("H7").Value = 126
$ws.Range("I7").Value = 162
$ws.Range("K7").Value = 18319
$ws.Range("L7").Value = 7042
$ws.Range("M7").Value = 11277
$ws.Range("N7").Value = 10629
$ws.Range("P7").Value = 583
$ws.Range("Q7").Value = 449
$ws.Range("R7").Value = -808
$ws.Range("S7").Value = 831
$ws.Range("T7").Value = 745
$ws.Range("U7").Value = -421
$ws.Range("W7").Value = 4.19
$ws.Range("X7").Value = 0.92
$ws.Range("Y7").Value = 1.52
$ws.Range("Z7").Value = 0.71
$ws.Range("AA7").Value = 62.45
$ws.Range("AC7").Value = 1383
$ws.Range("AD7").Value = 89.31
$ws.Range("AE7").Value = 93129
$ws.Range("AF7").Value = 1.33
$ws.Range("AG7").Value = 1060
$ws.Range("AH7").Value = 0.86
$ws.Range("AI7").Value = 76.66

$ws.Range("D8").Value = 14692
$ws.Range("E8").Value = 810
$ws.Range("G8").Value = 714
$ws.Range("H8").Value = 527
$ws.Range("I8").Value = 529
$ws.Range("K8").Value = 18731
$ws.Range("L8").Value = 7050
$ws.Range("M8").Value = 11679
$ws.Range("N8").Value = 11023
$ws.Range("P8").Value = 583
$ws.Range("Q8").Value = 478
$ws.Range("R8").Value = -420
$ws.Range("S8").Value = -169
$ws.Range("T8").Value = 321
$ws.Range("U8").Value = 93
$ws.Range("W8").Value = 5.51
$ws.Range("X8").Value = 3.59
$ws.Range("Y8").Value = 4.89
$ws.Range("Z8").Value = 2.84
$ws.Range("AA8").Value = 60.36
$ws.Range("AC8").Value = 4527
$ws.Range("AD8").Value = 27.28
$ws.Range("AE8").Value = 96584
$ws.Range("AF8").Value = 1.28
$ws.Range("AG8").Value = 1070
$ws.Range("AH8").Value = 0.87
$ws.Range("AI8").Value = 23.64

$ws.Range("D9").Value = 15683
$ws.Range("E9").Value = 887
$ws.Range("G9").Value = 763
$ws.Range("H9").Value = 564
$ws.Range("I9").Value = 579
$ws.Range("K9").Value = 19226
$ws.Range("L9").Value = 7117
$ws.Range("M9").Value = 12109
$ws.Range("N9").Value = 11456
$ws.Range("P9").Value = 583
$ws.Range("Q9").Value = 527
$ws.Range("R9").Value = -385
$ws.Range("S9").Value = -267
$ws.Range("T9").Value = 296
$ws.Range("U9").Value = 265
$ws.Range("W9").Value = 5.66
$ws.Range("X9").Value = 3.59
$ws.Range("Y9").Value = 5.15
$ws.Range("Z9").Value = 2.97
$ws.Range("AA9").Value = 58.78
$ws.Range("AC9").Value = 4951
$ws.Range("AD9").Value = 24.94
$ws.Range("AE9").Value = 100375
$ws.Range("AF9").Value = 1.23
$ws.Range("AG9").Value = 1110
$ws.Range("AH9").Value = 0.9
$ws.Range("AI9").Value = 22.42

# Clear cells removed in the target revision
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
